$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02743666666666666
$ws.Range("H2").Value = 0.08230999999999999
$ws.Range("I2").Value = 0.007366285056527356
$ws.Range("J2").Value = 0.007366285056527356
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.103724333333333
$ws.Range("N2").Value = 3.311173
$ws.Range("O2").Value = 0.01870879385910814
$ws.Range("P2").Value = 0.01870879385910814
$ws.Range("Q2").Value = 0.03028251662555555
$ws.Range("R2").Value = 0.27254264963
$ws.Range("S2").Value = 0.0001378143086299991
$ws.Range("T2").Value = 0.0001378143086299991
# Row 3
$ws.Range("G3").Value = 0.02743666666666666
$ws.Range("H3").Value = 0.08230999999999999
$ws.Range("I3").Value = 0.007366285056527356
$ws.Range("J3").Value = 0.007366285056527356
$ws.Range("O3").Value = 0.1603368629650925
$ws.Range("P3").Value = 0.1603368629650925
$ws.Range("Q3").Value = 0.2595252133833333
$ws.Range("R3").Value = 2.33572692045
$ws.Range("S3").Value = 0.001181087037670236
$ws.Range("T3").Value = 0.001181087037670236
# Row 4
$ws.Range("G4").Value = 0.02743666666666666
$ws.Range("H4").Value = 0.08230999999999999
$ws.Range("I4").Value = 0.007366285056527356
$ws.Range("J4").Value = 0.007366285056527356
$ws.Range("M4").Value = 47.61312599999999
$ws.Range("N4").Value = 142.839378
$ws.Range("O4").Value = 0.80707123365805
$ws.Range("P4").Value = 0.80707123365805
$ws.Range("Q4").Value = 1.30634546702
$ws.Range("R4").Value = 11.75710920318
$ws.Range("S4").Value = 0.005945116768048392
$ws.Range("T4").Value = 0.005945116768048392
# Row 5
$ws.Range("G5").Value = 0.02743666666666666
$ws.Range("H5").Value = 0.08230999999999999
$ws.Range("I5").Value = 0.007366285056527356
$ws.Range("J5").Value = 0.007366285056527356
$ws.Range("M5").Value = 0.8190333333333334
$ws.Range("N5").Value = 2.4571
$ws.Range("O5").Value = 0.01388310951774934
$ws.Range("P5").Value = 0.01388310951774934
$ws.Range("Q5").Value = 0.02247154455555556
$ws.Range("R5").Value = 0.202243901
$ws.Range("S5").Value = 0.0001022669421787296
$ws.Range("T5").Value = 0.0001022669421787296
# Row 6
$ws.Range("G6").Value = 3.368329
$ws.Range("I6").Value = 0.9043398704228307
$ws.Range("J6").Value = 0.9043398704228307
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.103724333333333
$ws.Range("N6").Value = 3.311173
$ws.Range("O6").Value = 0.01870879385910814
$ws.Range("P6").Value = 0.01870879385910814
$ws.Range("Q6").Value = 3.717706679972333
$ws.Range("R6").Value = 33.459360119751
$ws.Range("S6").Value = 0.01691910821431331
$ws.Range("T6").Value = 0.01691910821431331
# Row 7
$ws.Range("G7").Value = 3.368329
$ws.Range("I7").Value = 0.9043398704228307
$ws.Range("J7").Value = 0.9043398704228307
$ws.Range("O7").Value = 0.1603368629650925
$ws.Range("P7").Value = 0.1603368629650925
$ws.Range("Q7").Value = 31.861242952385
$ws.Range("S7").Value = 0.1449990178778549
$ws.Range("T7").Value = 0.1449990178778549
# Row 8
$ws.Range("G8").Value = 3.368329
$ws.Range("I8").Value = 0.9043398704228307
$ws.Range("J8").Value = 0.9043398704228307
$ws.Range("M8").Value = 47.61312599999999
$ws.Range("N8").Value = 142.839378
$ws.Range("O8").Value = 0.80707123365805
$ws.Range("P8").Value = 0.80707123365805
$ws.Range("Q8").Value = 160.376673086454
$ws.Range("R8").Value = 1443.390057778086
$ws.Range("S8").Value = 0.729866694868315
$ws.Range("T8").Value = 0.729866694868315
# Row 9
$ws.Range("G9").Value = 3.368329
$ws.Range("I9").Value = 0.9043398704228307
$ws.Range("J9").Value = 0.9043398704228307
$ws.Range("M9").Value = 0.8190333333333334
$ws.Range("N9").Value = 2.4571
$ws.Range("O9").Value = 0.01388310951774934
$ws.Range("P9").Value = 0.01388310951774934
$ws.Range("Q9").Value = 2.758773728633333
$ws.Range("R9").Value = 24.8289635577
$ws.Range("S9").Value = 0.0125550494623474
$ws.Range("T9").Value = 0.0125550494623474
# Row 10
$ws.Range("G10").Value = 0.3288616666666667
$ws.Range("H10").Value = 0.9865849999999999
$ws.Range("I10").Value = 0.08829384452064198
$ws.Range("J10").Value = 0.08829384452064198
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.103724333333333
$ws.Range("N10").Value = 3.311173
$ws.Range("O10").Value = 0.01870879385910814
$ws.Range("P10").Value = 0.01870879385910814
$ws.Range("Q10").Value = 0.3629726238005556
$ws.Range("R10").Value = 3.266753614205
$ws.Range("S10").Value = 0.001651871336164836
$ws.Range("T10").Value = 0.001651871336164836
# Row 11
$ws.Range("G11").Value = 0.3288616666666667
$ws.Range("H11").Value = 0.9865849999999999
$ws.Range("I11").Value = 0.08829384452064198
$ws.Range("J11").Value = 0.08829384452064198
$ws.Range("O11").Value = 0.1603368629650925
$ws.Range("P11").Value = 0.1603368629650925
$ws.Range("Q11").Value = 3.110723881008334
$ws.Range("R11").Value = 27.996514929075
$ws.Range("S11").Value = 0.01415675804956736
$ws.Range("T11").Value = 0.01415675804956736
# Row 12
$ws.Range("G12").Value = 0.3288616666666667
$ws.Range("H12").Value = 0.9865849999999999
$ws.Range("I12").Value = 0.08829384452064198
$ws.Range("J12").Value = 0.08829384452064198
$ws.Range("M12").Value = 47.61312599999999
$ws.Range("N12").Value = 142.839378
$ws.Range("O12").Value = 0.80707123365805
$ws.Range("P12").Value = 0.80707123365805
$ws.Range("Q12").Value = 15.65813197157
$ws.Range("R12").Value = 140.92318774413
$ws.Range("S12").Value = 0.07125942202168659
$ws.Range("T12").Value = 0.07125942202168659
# Row 13
$ws.Range("G13").Value = 0.3288616666666667
$ws.Range("H13").Value = 0.9865849999999999
$ws.Range("I13").Value = 0.08829384452064198
$ws.Range("J13").Value = 0.08829384452064198
$ws.Range("M13").Value = 0.8190333333333334
$ws.Range("N13").Value = 2.4571
$ws.Range("O13").Value = 0.01388310951774934
$ws.Range("P13").Value = 0.01388310951774934
$ws.Range("Q13").Value = 0.2693486670555555
$ws.Range("R13").Value = 2.4241380035
$ws.Range("S13").Value = 0.001225793113223205
$ws.Range("T13").Value = 0.001225793113223205
